# Update the metadata/date and the CodeSystem System URI values.
$wb = $excel.ActiveWorkbook

# Sheet1 "Metadata": row 8 holds "Date" (A8) / timestamp (B8)
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# Sheet2 "Include #0": row 4 System URI -> TRE-R21-Fonction
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R21-Fonction"

# Sheet3 "Include #1": row 4 System URI -> TRE-R96-AutreFonctionSanitaire
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R96-AutreFonctionSanitaire"

# Sheet4 "Include #2": row 4 System URI -> TRE-R85-RolePriseCharge
$wsInc2 = $wb.Worksheets.Item("Include #2")
$wsInc2.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R85-RolePriseCharge"
